$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.133054
$ws.Range("H2").Value = 0.399162
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.372179
$ws.Range("N2").Value = 4.116537
$ws.Range("O2").Value = 0.2533567233062949
$ws.Range("P2").Value = 0.2533567233062949
$ws.Range("Q2").Value = 0.182573904666
$ws.Range("R2").Value = 1.643165141994
$ws.Range("S2").Value = 0.2533567233062949
$ws.Range("T2").Value = 0.2533567233062949

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.133054
$ws.Range("H3").Value = 0.399162
$ws.Range("O3").Value = 0.01171683533985869
$ws.Range("P3").Value = 0.0117168353398587
$ws.Range("Q3").Value = 0.008443385083333333
$ws.Range("R3").Value = 0.07599046575
$ws.Range("S3").Value = 0.01171683533985869
$ws.Range("T3").Value = 0.0117168353398587

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.133054
$ws.Range("H4").Value = 0.399162
$ws.Range("M4").Value = 3.980358666666666
$ws.Range("O4").Value = 0.7349264413538463
$ws.Range("P4").Value = 0.7349264413538463
$ws.Range("Q4").Value = 0.5296026420346667
$ws.Range("R4").Value = 4.766423778311999
$ws.Range("S4").Value = 0.7349264413538463
$ws.Range("T4").Value = 0.7349264413538463
